$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title shape: merge the three runs "A" + " " + "slide" into a single run "A slide".
# Setting directly to the same concatenated text is a no-op for the run layout, so
# first shift through an unrelated placeholder value to force a single-run rewrite.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "x"
$titleRange.Text = "A slide"

# Table (shape 3), second cell of the first row: merge "a" + " " + "table" into "a table".
$cellRange = $s.Shapes.Item(3).Table.Cell(1, 2).Shape.TextFrame.TextRange
$cellRange.Text = "x"
$cellRange.Text = "a table"
